# Apply edits described by the commit diff for
# 2023/italy_serie-c-group-b_2023-2024.xlsx
#
# Summary of changes:
#  1. Three trios/pairs of rows had their match-data columns (F:V) rotated
#     or swapped between each other (rows 14/15/18, 67/68/69, 74/75,
#     80/81/82) while the leading Indice/pais/torneio/temporada/data_partida
#     columns (A:E) stayed untouched.
#  2. Three brand-new match rows (90, 91, 92) were appended at the bottom
#     of the sheet, extending the used range from A1:V89 to A1:V92.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows whose match content was reordered (rotations/swaps) ---

# Row 14
$ws.Cells.Item(14, 6).Value = "Torres"
$ws.Cells.Item(14, 7).Value = 2
$ws.Cells.Item(14, 8).Value = "Rimini"
$ws.Cells.Item(14, 9).Value = 1
$ws.Cells.Item(14, 10).Value = 1.98
$ws.Cells.Item(14, 11).Value = "07/09/2023 15:42"
$ws.Cells.Item(14, 12).Value = 1.77
$ws.Cells.Item(14, 13).Value = "09/09/2023 18:03"
$ws.Cells.Item(14, 14).Value = 3.04
$ws.Cells.Item(14, 15).Value = "07/09/2023 15:42"
$ws.Cells.Item(14, 16).Value = 3.26
$ws.Cells.Item(14, 17).Value = "09/09/2023 18:03"
$ws.Cells.Item(14, 18).Value = 3.88
$ws.Cells.Item(14, 19).Value = "07/09/2023 15:42"
$ws.Cells.Item(14, 20).Value = 5.4
$ws.Cells.Item(14, 21).Value = "09/09/2023 18:06"
$ws.Cells.Item(14, 22).Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/sassari-torres-rimini/CrECZhom/"

# Row 15
$ws.Cells.Item(15, 6).Value = "Fermana"
$ws.Cells.Item(15, 7).Value = 1
$ws.Cells.Item(15, 8).Value = "Pontedera"
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 4.23
$ws.Cells.Item(15, 11).Value = "08/09/2023 06:42"
$ws.Cells.Item(15, 12).Value = 4.21
$ws.Cells.Item(15, 13).Value = "09/09/2023 18:29"
$ws.Cells.Item(15, 14).Value = 3.02
$ws.Cells.Item(15, 15).Value = "08/09/2023 06:42"
$ws.Cells.Item(15, 16).Value = 3.2
$ws.Cells.Item(15, 17).Value = "09/09/2023 18:29"
$ws.Cells.Item(15, 18).Value = 1.9
$ws.Cells.Item(15, 19).Value = "08/09/2023 06:42"
$ws.Cells.Item(15, 20).Value = 1.97
$ws.Cells.Item(15, 21).Value = "09/09/2023 18:29"
$ws.Cells.Item(15, 22).Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/fermana-us-pontedera/S2jkVwf1/"

# Row 18
$ws.Cells.Item(18, 6).Value = "Sestri Levante"
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = "Lucchese"
$ws.Cells.Item(18, 9).Value = 1
$ws.Cells.Item(18, 10).Value = 3.29
$ws.Cells.Item(18, 11).Value = "07/09/2023 14:42"
$ws.Cells.Item(18, 12).Value = 3.54
$ws.Cells.Item(18, 13).Value = "09/09/2023 18:06"
$ws.Cells.Item(18, 14).Value = 2.94
$ws.Cells.Item(18, 15).Value = "07/09/2023 14:42"
$ws.Cells.Item(18, 16).Value = 3.17
$ws.Cells.Item(18, 17).Value = "09/09/2023 18:06"
$ws.Cells.Item(18, 18).Value = 2.23
$ws.Cells.Item(18, 19).Value = "07/09/2023 14:42"
$ws.Cells.Item(18, 20).Value = 2.18
$ws.Cells.Item(18, 21).Value = "09/09/2023 18:06"
$ws.Cells.Item(18, 22).Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/sestri-levante-lucchese/Ece9zZVt/"

# Row 67
$ws.Cells.Item(67, 6).Value = "Gubbio"
$ws.Cells.Item(67, 7).Value = 2
$ws.Cells.Item(67, 8).Value = "Carrarese"
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 2.49
$ws.Cells.Item(67, 11).Value = "05/10/2023 15:12"
$ws.Cells.Item(67, 12).Value = 2.32
$ws.Cells.Item(67, 13).Value = "09/10/2023 19:22"
$ws.Cells.Item(67, 14).Value = 2.83
$ws.Cells.Item(67, 15).Value = "05/10/2023 15:12"
$ws.Cells.Item(67, 16).Value = 2.9
$ws.Cells.Item(67, 17).Value = "09/10/2023 20:31"
$ws.Cells.Item(67, 18).Value = 2.88
$ws.Cells.Item(67, 19).Value = "05/10/2023 15:12"
$ws.Cells.Item(67, 20).Value = 3.54
$ws.Cells.Item(67, 21).Value = "09/10/2023 19:22"
$ws.Cells.Item(67, 22).Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/gubbio-carrarese/U1wX2FB1/"

# Row 68
$ws.Cells.Item(68, 6).Value = "Arezzo"
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = "Cesena"
$ws.Cells.Item(68, 9).Value = 2
$ws.Cells.Item(68, 10).Value = 4.99
$ws.Cells.Item(68, 11).Value = "05/10/2023 15:12"
$ws.Cells.Item(68, 12).Value = 6.14
$ws.Cells.Item(68, 13).Value = "09/10/2023 20:43"
$ws.Cells.Item(68, 14).Value = 3.49
$ws.Cells.Item(68, 15).Value = "05/10/2023 15:12"
$ws.Cells.Item(68, 16).Value = 4.21
$ws.Cells.Item(68, 17).Value = "09/10/2023 20:43"
$ws.Cells.Item(68, 18).Value = 1.65
$ws.Cells.Item(68, 19).Value = "05/10/2023 15:12"
$ws.Cells.Item(68, 20).Value = 1.52
$ws.Cells.Item(68, 21).Value = "09/10/2023 20:43"
$ws.Cells.Item(68, 22).Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/arezzo-cesena/AaNV4ysk/"

# Row 69
$ws.Cells.Item(69, 6).Value = "Fermana"
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = "Perugia"
$ws.Cells.Item(69, 9).Value = 2
$ws.Cells.Item(69, 10).Value = 4.83
$ws.Cells.Item(69, 11).Value = "05/10/2023 15:12"
$ws.Cells.Item(69, 12).Value = 6.25
$ws.Cells.Item(69, 13).Value = "09/10/2023 20:44"
$ws.Cells.Item(69, 14).Value = 3.45
$ws.Cells.Item(69, 15).Value = "05/10/2023 15:12"
$ws.Cells.Item(69, 16).Value = 3.7
$ws.Cells.Item(69, 17).Value = "09/10/2023 20:44"
$ws.Cells.Item(69, 18).Value = 1.65
$ws.Cells.Item(69, 19).Value = "05/10/2023 15:12"
$ws.Cells.Item(69, 20).Value = 1.6
$ws.Cells.Item(69, 21).Value = "09/10/2023 20:28"
$ws.Cells.Item(69, 22).Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/fermana-perugia/ljMZ3ede/"

# Row 74
$ws.Cells.Item(74, 6).Value = "Carrarese"
$ws.Cells.Item(74, 7).Value = 1
$ws.Cells.Item(74, 8).Value = "Ancona"
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 1.74
$ws.Cells.Item(74, 11).Value = "12/10/2023 08:13"
$ws.Cells.Item(74, 12).Value = 1.72
$ws.Cells.Item(74, 13).Value = "15/10/2023 16:06"
$ws.Cells.Item(74, 14).Value = 3.32
$ws.Cells.Item(74, 15).Value = "12/10/2023 08:13"
$ws.Cells.Item(74, 16).Value = 3.57
$ws.Cells.Item(74, 17).Value = "15/10/2023 16:06"
$ws.Cells.Item(74, 18).Value = 4.46
$ws.Cells.Item(74, 19).Value = "12/10/2023 08:13"
$ws.Cells.Item(74, 20).Value = 5.14
$ws.Cells.Item(74, 21).Value = "15/10/2023 16:06"
$ws.Cells.Item(74, 22).Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/carrarese-ancona/Sd91IfBl/"

# Row 75
$ws.Cells.Item(75, 6).Value = "Recanatese"
$ws.Cells.Item(75, 7).Value = 2
$ws.Cells.Item(75, 8).Value = "Arezzo"
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 2.24
$ws.Cells.Item(75, 11).Value = "12/10/2023 08:12"
$ws.Cells.Item(75, 12).Value = 2.6
$ws.Cells.Item(75, 13).Value = "15/10/2023 16:14"
$ws.Cells.Item(75, 14).Value = 2.95
$ws.Cells.Item(75, 15).Value = "12/10/2023 08:12"
$ws.Cells.Item(75, 16).Value = 3.12
$ws.Cells.Item(75, 17).Value = "15/10/2023 16:14"
$ws.Cells.Item(75, 18).Value = 3.13
$ws.Cells.Item(75, 19).Value = "12/10/2023 08:12"
$ws.Cells.Item(75, 20).Value = 2.84
$ws.Cells.Item(75, 21).Value = "15/10/2023 16:14"
$ws.Cells.Item(75, 22).Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/recanatese-arezzo/jkc8yzIQ/"

# Row 80
$ws.Cells.Item(80, 6).Value = "Fermana"
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = "Entella"
$ws.Cells.Item(80, 9).Value = 2
$ws.Cells.Item(80, 10).Value = 5.14
$ws.Cells.Item(80, 11).Value = "19/10/2023 08:13"
$ws.Cells.Item(80, 12).Value = 5
$ws.Cells.Item(80, 13).Value = "22/10/2023 13:52"
$ws.Cells.Item(80, 14).Value = 3.29
$ws.Cells.Item(80, 15).Value = "19/10/2023 08:13"
$ws.Cells.Item(80, 16).Value = 3.27
$ws.Cells.Item(80, 17).Value = "22/10/2023 13:52"
$ws.Cells.Item(80, 18).Value = 1.68
$ws.Cells.Item(80, 19).Value = "19/10/2023 08:13"
$ws.Cells.Item(80, 20).Value = 1.81
$ws.Cells.Item(80, 21).Value = "22/10/2023 13:52"
$ws.Cells.Item(80, 22).Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/fermana-entella/djxmUiA0/"

# Row 81
$ws.Cells.Item(81, 6).Value = "Juventus U23"
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = "Perugia"
$ws.Cells.Item(81, 9).Value = 2
$ws.Cells.Item(81, 10).Value = 3.16
$ws.Cells.Item(81, 11).Value = "19/10/2023 08:13"
$ws.Cells.Item(81, 12).Value = 3.85
$ws.Cells.Item(81, 13).Value = "22/10/2023 13:54"
$ws.Cells.Item(81, 14).Value = 2.98
$ws.Cells.Item(81, 15).Value = "19/10/2023 08:13"
$ws.Cells.Item(81, 16).Value = 3.36
$ws.Cells.Item(81, 17).Value = "22/10/2023 13:54"
$ws.Cells.Item(81, 18).Value = 2.22
$ws.Cells.Item(81, 19).Value = "19/10/2023 08:13"
$ws.Cells.Item(81, 20).Value = 2
$ws.Cells.Item(81, 21).Value = "22/10/2023 13:54"
$ws.Cells.Item(81, 22).Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/juventus-perugia/I5QeSVvD/"

# Row 82
$ws.Cells.Item(82, 6).Value = "Torres"
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = "Pontedera"
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 10).Value = 1.84
$ws.Cells.Item(82, 11).Value = "19/10/2023 08:13"
$ws.Cells.Item(82, 12).Value = 1.99
$ws.Cells.Item(82, 13).Value = "22/10/2023 13:53"
$ws.Cells.Item(82, 14).Value = 3.18
$ws.Cells.Item(82, 15).Value = "19/10/2023 08:13"
$ws.Cells.Item(82, 16).Value = 3.12
$ws.Cells.Item(82, 17).Value = "22/10/2023 13:53"
$ws.Cells.Item(82, 18).Value = 4.12
$ws.Cells.Item(82, 19).Value = "19/10/2023 08:13"
$ws.Cells.Item(82, 20).Value = 4.29
$ws.Cells.Item(82, 21).Value = "22/10/2023 13:53"
$ws.Cells.Item(82, 22).Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/sassari-torres-us-pontedera/rmbIqX9m/"

# --- Add new rows 90-92 (new matches appended) ---
$ws.Range("A89:V89").Copy($ws.Range("A90:V90"))
$ws.Range("A89:V89").Copy($ws.Range("A91:V91"))
$ws.Range("A89:V89").Copy($ws.Range("A92:V92"))

# Row 90
$ws.Cells.Item(90, 1).Value = 89
$ws.Cells.Item(90, 2).Value = "italy"
$ws.Cells.Item(90, 3).Value = "serie-c-group-b"
$ws.Cells.Item(90, 4).Value = "2023-2024"
$ws.Cells.Item(90, 5).Value = 45224.77083333334
$ws.Cells.Item(90, 6).Value = "Gubbio"
$ws.Cells.Item(90, 7).Value = 1
$ws.Cells.Item(90, 8).Value = "Juventus U23"
$ws.Cells.Item(90, 9).Value = 1
$ws.Cells.Item(90, 10).Value = 1.93
$ws.Cells.Item(90, 11).Value = "24/10/2023 12:42"
$ws.Cells.Item(90, 12).Value = 1.72
$ws.Cells.Item(90, 13).Value = "25/10/2023 18:24"
$ws.Cells.Item(90, 14).Value = 3.1
$ws.Cells.Item(90, 15).Value = "24/10/2023 12:42"
$ws.Cells.Item(90, 16).Value = 3.56
$ws.Cells.Item(90, 17).Value = "25/10/2023 18:24"
$ws.Cells.Item(90, 18).Value = 4.02
$ws.Cells.Item(90, 19).Value = "24/10/2023 12:42"
$ws.Cells.Item(90, 20).Value = 5.07
$ws.Cells.Item(90, 21).Value = "25/10/2023 18:24"
$ws.Cells.Item(90, 22).Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/gubbio-juventus/AakZuk9C/"

# Row 91
$ws.Cells.Item(91, 1).Value = 90
$ws.Cells.Item(91, 2).Value = "italy"
$ws.Cells.Item(91, 3).Value = "serie-c-group-b"
$ws.Cells.Item(91, 4).Value = "2023-2024"
$ws.Cells.Item(91, 5).Value = 45224.77083333334
$ws.Cells.Item(91, 6).Value = "Olbia"
$ws.Cells.Item(91, 7).Value = 2
$ws.Cells.Item(91, 8).Value = "Fermana"
$ws.Cells.Item(91, 9).Value = 1
$ws.Cells.Item(91, 10).Value = 2.07
$ws.Cells.Item(91, 11).Value = "24/10/2023 12:42"
$ws.Cells.Item(91, 12).Value = 2.1
$ws.Cells.Item(91, 13).Value = "25/10/2023 18:27"
$ws.Cells.Item(91, 14).Value = 2.96
$ws.Cells.Item(91, 15).Value = "24/10/2023 12:42"
$ws.Cells.Item(91, 16).Value = 2.95
$ws.Cells.Item(91, 17).Value = "25/10/2023 18:28"
$ws.Cells.Item(91, 18).Value = 3.55
$ws.Cells.Item(91, 19).Value = "24/10/2023 12:42"
$ws.Cells.Item(91, 20).Value = 4.14
$ws.Cells.Item(91, 21).Value = "25/10/2023 18:28"
$ws.Cells.Item(91, 22).Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/olbia-fermana/vZgwu9OI/"

# Row 92
$ws.Cells.Item(92, 1).Value = 91
$ws.Cells.Item(92, 2).Value = "italy"
$ws.Cells.Item(92, 3).Value = "serie-c-group-b"
$ws.Cells.Item(92, 4).Value = "2023-2024"
$ws.Cells.Item(92, 5).Value = 45224.86458333334
$ws.Cells.Item(92, 6).Value = "Pescara"
$ws.Cells.Item(92, 7).Value = 1
$ws.Cells.Item(92, 8).Value = "Torres"
$ws.Cells.Item(92, 9).Value = 2
$ws.Cells.Item(92, 10).Value = 1.88
$ws.Cells.Item(92, 11).Value = "24/10/2023 12:42"
$ws.Cells.Item(92, 12).Value = 2.05
$ws.Cells.Item(92, 13).Value = "25/10/2023 20:37"
$ws.Cells.Item(92, 14).Value = 3.39
$ws.Cells.Item(92, 15).Value = "24/10/2023 12:42"
$ws.Cells.Item(92, 16).Value = 3.4
$ws.Cells.Item(92, 17).Value = "25/10/2023 20:42"
$ws.Cells.Item(92, 18).Value = 3.61
$ws.Cells.Item(92, 19).Value = "24/10/2023 12:42"
$ws.Cells.Item(92, 20).Value = 3.62
$ws.Cells.Item(92, 21).Value = "25/10/2023 20:37"
$ws.Cells.Item(92, 22).Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/pescara-sassari-torres/IPhsvTwP/"

